# Applies the "adds info to project metadata" commit to the feather metadata workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. personnel sheet: role changes from "project lead" to "creator"
# ---------------------------------------------------------------------------
$wsPersonnel = $wb.Worksheets.Item("personnel")
$wsPersonnel.Range("D2").Value = "creator"

# ---------------------------------------------------------------------------
# 2. project sheet: role changes from "Primary Investigator" to "creator"
# ---------------------------------------------------------------------------
$wsProject = $wb.Worksheets.Item("project")
$wsProject.Range("D2").Value = "creator"
$wsProject.Range("F2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. maintenance sheet: add status value "complete"
# ---------------------------------------------------------------------------
$wsMaintenance = $wb.Worksheets.Item("maintenance")
$wsMaintenance.Range("A2").Value = "complete"
$wsMaintenance.Range("B2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. coverage sheet: add geographic description + bounding coordinates
# ---------------------------------------------------------------------------
$wsCoverage = $wb.Worksheets.Item("coverage")
$wsCoverage.Range("A2").Value = "Feather river "
$wsCoverage.Range("B2").Value = -121.63263600000001
$wsCoverage.Range("C2").Value = -121.60463300000001
$wsCoverage.Range("D2").Value = 39.4621
$wsCoverage.Range("E2").Value = 39.212150000000001
$wsCoverage.Range("A3").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. title sheet: selection bookkeeping
# ---------------------------------------------------------------------------
$wsTitle = $wb.Worksheets.Item("title")
$wsTitle.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 6. personnel becomes the active/selected sheet (matches activeTab="1")
# ---------------------------------------------------------------------------
$wsPersonnel.Activate() | Out-Null
$wsPersonnel.Range("G2").Select() | Out-Null
